$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("A1").Value = "material_name"
$ws.Range("B1").Value = "specification"
$ws.Range("C1").Value = "quantity"
$ws.Range("D1").Value = "unit"

# Row 2: shift data - A2 becomes old B2 text, B2 becomes old C2 text, C2 becomes old E2 value, D2 stays
$ws.Range("A2").Value = "卡箍"
$ws.Range("B2").Value = "DN100"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = "个"

# Row 3
$ws.Range("A3").Value = "沟槽大小头"
$ws.Range("B3").Value = "DN100*80"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "个"

# Row 4
$ws.Range("A4").Value = "沟槽弯头"
$ws.Range("B4").Value = "DN80*65"
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = "个"

# Clear columns E and F entirely (old "单位", "数量", "连接方式" data)
$ws.Range("E1:F4").Clear()
